$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 54.75
$ws.Range("I4").Value = 54.75
$ws.Range("K4").Value = 54.75
$ws.Range("M4").Value = 59.25
$ws.Range("H17").Value = 1166.2858
$ws.Range("J17").Value = 1166.2858
$ws.Range("L17").Value = 3498.8574
$ws.Range("N17").Value = -3834.8574
$ws.Range("H43").Value = 1927200.9
$ws.Range("I43").Value = 3849403.2
$ws.Range("J43").Value = 4998.5
$ws.Range("K43").Value = 3849403.2
$ws.Range("L43").Value = 4998.5
$ws.Range("M43").Value = -3849334.2
$ws.Range("N43").Value = -5136.5
$ws.Range("H80").Value = 794.13336
$ws.Range("I80").Value = 730.3333
$ws.Range("J80").Value = 889.8333
$ws.Range("K80").Value = 2190.9999
$ws.Range("L80").Value = 2669.4999
$ws.Range("M80").Value = -1192.9999
$ws.Range("N80").Value = -4665.4999
$ws.Range("H83").Value = 794.13336
$ws.Range("I83").Value = 730.3333
$ws.Range("J83").Value = 889.8333
$ws.Range("K83").Value = 6572.9997
$ws.Range("L83").Value = 8008.4997
$ws.Range("M83").Value = -1580.9997
$ws.Range("N83").Value = -17992.4997
$ws.Range("H94").Value = 27787276
$ws.Range("I94").Value = 55557056
$ws.Range("J94").Value = 17497
$ws.Range("K94").Value = 55557056
$ws.Range("L94").Value = 17497
$ws.Range("M94").Value = -55556605
$ws.Range("N94").Value = -18399
$ws.Range("H98").Value = 1107.6428
$ws.Range("I98").Value = 864.36365
$ws.Range("K98").Value = 864.36365
$ws.Range("M98").Value = 633.63635
$ws.Range("H116").Value = 34639544
$ws.Range("I116").Value = 28971468
$ws.Range("K116").Value = 28971468
$ws.Range("M116").Value = -28968026
$ws.Range("H122").Value = 1107.6428
$ws.Range("I122").Value = 864.36365
$ws.Range("K122").Value = 2593.09095
$ws.Range("M122").Value = -143.0909499999998
$ws.Range("H132").Value = 14980.64
$ws.Range("I132").Value = 2992.5652
$ws.Range("K132").Value = 8977.695599999999
$ws.Range("M132").Value = -6447.695599999999
$ws.Range("H136").Value = 129550
$ws.Range("J136").Value = 129550
$ws.Range("L136").Value = 129550
$ws.Range("N136").Value = -139750
$ws.Range("H140").Value = 66701
$ws.Range("J140").Value = 66128.42999999999
$ws.Range("L140").Value = 66128.42999999999
$ws.Range("N140").Value = -76488.42999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6010.4375
$ws.Range("J32").Value = 10008.866
$ws.Range("L32").Value = 10008.866
$ws.Range("N32").Value = -10582.866
$ws.Range("H61").Value = 11413.533
$ws.Range("I61").Value = 12086
$ws.Range("K61").Value = 12086
$ws.Range("M61").Value = -11874
$ws.Range("H63").Value = 3399.5
$ws.Range("I63").Value = 3399.5
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3399.5
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2713.5
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 3399.5
$ws.Range("I66").Value = 3399.5
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 16997.5
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -13565.5
$ws.Range("N66").ClearContents()
$ws.Range("H132").Value = 3455.6128
$ws.Range("I132").Value = 2033.0625
$ws.Range("K132").Value = 6099.1875
$ws.Range("M132").Value = -3569.1875
$ws.Range("H136").Value = 11413.533
$ws.Range("I136").Value = 12086
$ws.Range("K136").Value = 36258
$ws.Range("M136").Value = -33708
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 75855.8
$ws.Range("J60").Value = 82319.75
$ws.Range("L60").Value = 82319.75
$ws.Range("N60").Value = -83517.75
$ws.Range("H86").Value = 1483.3077
$ws.Range("J86").Value = 1099
$ws.Range("L86").Value = 1099
$ws.Range("N86").Value = -3345
$ws.Range("H89").Value = 1483.3077
$ws.Range("J89").Value = 1099
$ws.Range("L89").Value = 5495
$ws.Range("N89").Value = -16727
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24396644
$ws.Range("I31").Value = 90910600
$ws.Range("J31").Value = 8193.6
$ws.Range("K31").Value = 90910600
$ws.Range("L31").Value = 8193.6
$ws.Range("M31").Value = -90910305
$ws.Range("N31").Value = -8783.6
$ws.Range("H34").Value = 24396644
$ws.Range("I34").Value = 90910600
$ws.Range("J34").Value = 8193.6
$ws.Range("K34").Value = 90910600
$ws.Range("L34").Value = 8193.6
$ws.Range("M34").Value = -90910398
$ws.Range("N34").Value = -8597.6
$ws.Range("H58").Value = 630154.75
$ws.Range("J58").Value = 5498.4
$ws.Range("L58").Value = 5498.4
$ws.Range("N58").Value = -5904.4
$ws.Range("H99").Value = 6356.75
$ws.Range("I99").Value = 3570.8
$ws.Range("K99").Value = 3570.8
$ws.Range("M99").Value = -2072.8
$ws.Range("H126").Value = 6356.75
$ws.Range("I126").Value = 3570.8
$ws.Range("K126").Value = 10712.4
$ws.Range("M126").Value = -8242.400000000001
$ws.Range("H134").Value = 3950.9
$ws.Range("I134").Value = 4089.889
$ws.Range("J134").Value = 2700
$ws.Range("K134").Value = 12269.667
$ws.Range("L134").Value = 8100
$ws.Range("M134").Value = -9734.667000000001
$ws.Range("N134").Value = -13170
$ws.Range("H136").Value = 630154.75
$ws.Range("J136").Value = 5498.4
$ws.Range("L136").Value = 16495.2
$ws.Range("N136").Value = -21595.2
$ws.Range("H141").Value = 113095.45
$ws.Range("J141").Value = 113095.45
$ws.Range("L141").Value = 113095.45
$ws.Range("N141").Value = -123455.45
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 482.375
$ws.Range("J9").Value = 547.3333
$ws.Range("L9").Value = 1641.9999
$ws.Range("N9").Value = -2089.9999
$ws.Range("H140").Value = 3605.4324
$ws.Range("I140").Value = 2641.1738
$ws.Range("J140").Value = 5189.5713
$ws.Range("K140").Value = 7923.5214
$ws.Range("L140").Value = 15568.7139
$ws.Range("M140").Value = -2743.5214
$ws.Range("N140").Value = -25928.7139
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1768786
$ws.Range("I70").Value = 2511143.5
$ws.Range("K70").Value = 2511143.5
$ws.Range("M70").Value = -2510873.5
$ws.Range("H73").Value = 1768786
$ws.Range("I73").Value = 2511143.5
$ws.Range("K73").Value = 2511143.5
$ws.Range("M73").Value = -2510207.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1119.25
$ws.Range("I22").Value = 999.8461
$ws.Range("J22").Value = 1341
$ws.Range("K22").Value = 999.8461
$ws.Range("L22").Value = 1341
$ws.Range("M22").Value = -704.8461
$ws.Range("N22").Value = -1931
$ws.Range("H27").Value = 1119.25
$ws.Range("I27").Value = 999.8461
$ws.Range("J27").Value = 1341
$ws.Range("K27").Value = 999.8461
$ws.Range("L27").Value = 1341
$ws.Range("M27").Value = -892.8461
$ws.Range("N27").Value = -1555
$ws.Range("H93").Value = 1518.1666
$ws.Range("I93").Value = 1276.8
$ws.Range("K93").Value = 1276.8
$ws.Range("M93").Value = -28.79999999999995
